# Auto-generated PowerShell COM-interop script
# Applies "Stricter separation of departement specific pre mid and post mid courses"
$wb = $excel.ActiveWorkbook

# ---- Section_A: timetable classroom reassignments ----
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("E3").Value = "CS161 [C202]"
$ws.Range("D4").Value = "EC161 [C003]"
$ws.Range("C6").Value = "MA162 [C003]"
$ws.Range("E6").Value = "CS161 (Lab) [L106]"
$ws.Range("E7").Value = "CS161 (Lab) [L106]"
$ws.Range("B8").Value = "MA162 [C003]"
$ws.Range("C8").Value = "EC161 [C003]"
$ws.Range("D8").Value = "CS161 [C202]"

# ---- Section_B: timetable classroom reassignments ----
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("C3").Value = "MA161 [C003]"
$ws.Range("E3").Value = "CS161 [C203]"
$ws.Range("B4").Value = "MA161 [C003]"
$ws.Range("C4").Value = "DS161 [C003]"
$ws.Range("D4").Value = "EC161 [C004]"
$ws.Range("B6").Value = "DS161 [C003]"
$ws.Range("C6").Value = "MA162 [C004]"
$ws.Range("E6").Value = "CS161 (Lab) [L207]"
$ws.Range("E7").Value = "CS161 (Lab) [L207]"
$ws.Range("B8").Value = "MA162 [C004]"
$ws.Range("C8").Value = "EC161 [C004]"
$ws.Range("D8").Value = "CS161 [C203]"

# ---- Classroom_Utilization: recomputed usage stats for reassigned rooms ----
$ws = $wb.Worksheets.Item("Classroom_Utilization")
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 2.4
$ws.Range("G4").Value = 30
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 2.4
$ws.Range("G5").Value = 30
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0.6
$ws.Range("G14").Value = 7.5
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 0.6
$ws.Range("G15").Value = 7.5
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0.6
$ws.Range("G21").Value = 7.5

# ---- Classroom_Allocation: per-session room/facility reassignments ----
$ws = $wb.Worksheets.Item("Classroom_Allocation")
$ws.Range("G4").Value = "large classroom"
$ws.Range("H4").Value = "135"
$ws.Range("I4").Value = "Projector"
$ws.Range("M4").Value = "C003"
$ws.Range("G7").Value = "large classroom"
$ws.Range("H7").Value = "135"
$ws.Range("I7").Value = "Projector"
$ws.Range("M7").Value = "C003"
$ws.Range("G8").Value = "large classroom"
$ws.Range("H8").Value = "135"
$ws.Range("I8").Value = "Projector"
$ws.Range("M8").Value = "C003"
$ws.Range("G9").Value = "large classroom"
$ws.Range("H9").Value = "135"
$ws.Range("I9").Value = "Projector"
$ws.Range("M9").Value = "C003"
$ws.Range("G12").Value = "classroom"
$ws.Range("H12").Value = "96"
$ws.Range("I12").Value = "Projector"
$ws.Range("M12").Value = "C202"
$ws.Range("G13").Value = "classroom"
$ws.Range("H13").Value = "96"
$ws.Range("I13").Value = "Projector"
$ws.Range("M13").Value = "C202"
$ws.Range("M14").Value = "L106"
$ws.Range("M15").Value = "L106"
$ws.Range("G16").Value = "Auditorium"
$ws.Range("H16").Value = "240"
$ws.Range("I16").Value = "Audio/Video System"
$ws.Range("M16").Value = "C004"
$ws.Range("I17").Value = $null
$ws.Range("M17").Value = "C001"
$ws.Range("G18").Value = "large classroom"
$ws.Range("H18").Value = "120"
$ws.Range("M18").Value = "C002"
$ws.Range("M19").Value = "C101"
$ws.Range("M20").Value = "C102"
$ws.Range("M21").Value = "C104"
$ws.Range("I22").Value = "Projector"
$ws.Range("M22").Value = "C202"
$ws.Range("M23").Value = "C203"
$ws.Range("M24").Value = "C204"
$ws.Range("H25").Value = "96"
$ws.Range("M25").Value = "C205"
$ws.Range("G26").Value = "Auditorium"
$ws.Range("H26").Value = "240"
$ws.Range("I26").Value = "Audio/Video System"
$ws.Range("M26").Value = "C004"
$ws.Range("I27").Value = $null
$ws.Range("M27").Value = "C001"
$ws.Range("G28").Value = "large classroom"
$ws.Range("H28").Value = "120"
$ws.Range("M28").Value = "C002"
$ws.Range("M29").Value = "C101"
$ws.Range("M30").Value = "C102"
$ws.Range("M31").Value = "C104"
$ws.Range("I32").Value = "Projector"
$ws.Range("M32").Value = "C202"
$ws.Range("M33").Value = "C203"
$ws.Range("M34").Value = "C204"
$ws.Range("H35").Value = "96"
$ws.Range("M35").Value = "C205"
$ws.Range("G36").Value = "Auditorium"
$ws.Range("H36").Value = "240"
$ws.Range("I36").Value = "Audio/Video System"
$ws.Range("M36").Value = "C004"
$ws.Range("I37").Value = $null
$ws.Range("M37").Value = "C001"
$ws.Range("G38").Value = "large classroom"
$ws.Range("H38").Value = "120"
$ws.Range("M38").Value = "C002"
$ws.Range("M39").Value = "C101"
$ws.Range("M40").Value = "C102"
$ws.Range("M41").Value = "C104"
$ws.Range("I42").Value = "Projector"
$ws.Range("M42").Value = "C202"
$ws.Range("M43").Value = "C203"
$ws.Range("M44").Value = "C204"
$ws.Range("H45").Value = "96"
$ws.Range("M45").Value = "C205"
$ws.Range("H46").Value = "135"
$ws.Range("I46").Value = "Projector"
$ws.Range("M46").Value = "C003"
$ws.Range("H47").Value = "135"
$ws.Range("I47").Value = "Projector"
$ws.Range("M47").Value = "C003"
$ws.Range("G48").Value = "Auditorium"
$ws.Range("H48").Value = "240"
$ws.Range("I48").Value = "Audio/Video System"
$ws.Range("M48").Value = "C004"
$ws.Range("H49").Value = "135"
$ws.Range("I49").Value = "Projector"
$ws.Range("M49").Value = "C003"
$ws.Range("H50").Value = "135"
$ws.Range("I50").Value = "Projector"
$ws.Range("M50").Value = "C003"
$ws.Range("G51").Value = "Auditorium"
$ws.Range("H51").Value = "240"
$ws.Range("I51").Value = "Audio/Video System"
$ws.Range("M51").Value = "C004"
$ws.Range("G52").Value = "Auditorium"
$ws.Range("H52").Value = "240"
$ws.Range("I52").Value = "Audio/Video System"
$ws.Range("M52").Value = "C004"
$ws.Range("G53").Value = "Auditorium"
$ws.Range("H53").Value = "240"
$ws.Range("I53").Value = "Audio/Video System"
$ws.Range("M53").Value = "C004"
$ws.Range("G56").Value = "classroom"
$ws.Range("H56").Value = "96"
$ws.Range("I56").Value = "TV"
$ws.Range("M56").Value = "C203"
$ws.Range("G57").Value = "classroom"
$ws.Range("H57").Value = "96"
$ws.Range("I57").Value = "TV"
$ws.Range("M57").Value = "C203"
$ws.Range("M58").Value = "L207"
$ws.Range("M59").Value = "L207"
$ws.Range("G60").Value = "classroom"
$ws.Range("H60").Value = "96"
$ws.Range("I60").Value = "Projector"
$ws.Range("M60").Value = "C302"
$ws.Range("I61").Value = "TV"
$ws.Range("M61").Value = "C303"
$ws.Range("M62").Value = "C304"
$ws.Range("M63").Value = "C305"
$ws.Range("G64").Value = "large classroom"
$ws.Range("H64").Value = "135"
$ws.Range("I64").Value = "Projector"
$ws.Range("M64").Value = "C003"
$ws.Range("G65").Value = "classroom"
$ws.Range("H65").Value = "80"
$ws.Range("I65").Value = "TV"
$ws.Range("M65").Value = "L402"
$ws.Range("M66").Value = "L403"
$ws.Range("G70").Value = "classroom"
$ws.Range("H70").Value = "96"
$ws.Range("I70").Value = "Projector"
$ws.Range("M70").Value = "C302"
$ws.Range("I71").Value = "TV"
$ws.Range("M71").Value = "C303"
$ws.Range("M72").Value = "C304"
$ws.Range("M73").Value = "C305"
$ws.Range("G74").Value = "large classroom"
$ws.Range("H74").Value = "135"
$ws.Range("I74").Value = "Projector"
$ws.Range("M74").Value = "C003"
$ws.Range("G75").Value = "classroom"
$ws.Range("H75").Value = "80"
$ws.Range("I75").Value = "TV"
$ws.Range("M75").Value = "L402"
$ws.Range("M76").Value = "L403"
$ws.Range("G80").Value = "classroom"
$ws.Range("H80").Value = "96"
$ws.Range("I80").Value = "Projector"
$ws.Range("M80").Value = "C302"
$ws.Range("I81").Value = "TV"
$ws.Range("M81").Value = "C303"
$ws.Range("M82").Value = "C304"
$ws.Range("M83").Value = "C305"
$ws.Range("G84").Value = "large classroom"
$ws.Range("H84").Value = "135"
$ws.Range("I84").Value = "Projector"
$ws.Range("M84").Value = "C003"
$ws.Range("G85").Value = "classroom"
$ws.Range("H85").Value = "80"
$ws.Range("I85").Value = "TV"
$ws.Range("M85").Value = "L402"
$ws.Range("M86").Value = "L403"

# ---- Basket_Course_Allocations: room pairs per elective course ----
$ws = $wb.Worksheets.Item("Basket_Course_Allocations")
$ws.Range("C2").Value = "C004, C302"
$ws.Range("C3").Value = "C001, C303"
$ws.Range("C4").Value = "C002, C304"
$ws.Range("C5").Value = "C101, C305"
$ws.Range("C6").Value = "C003, C102"
$ws.Range("C7").Value = "C104, L402"
$ws.Range("C8").Value = "C202, L403"
$ws.Range("C9").Value = "C203, L404"
$ws.Range("C10").Value = "C204, L405"
$ws.Range("C11").Value = "C205, L406"

Write-Host "Applied all timetable/classroom reassignment edits"